# Update cryptos list (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'65.768.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.03%  "
$ws.Range("E2").Style = "Normal"

# Row 3 - Ethereum
$ws.Range("D3").Value = "'2.679.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.54%  "
$ws.Range("E3").Style = "Normal"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"

# Row 5 - BNB
$ws.Range("D5").Value = "'600.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.89%  "
$ws.Range("E5").Style = "Normal"

# Row 6 - Solana
$ws.Range("D6").Value = "'156.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.48%  "
$ws.Range("E6").Style = "Normal"

# Row 7 - USDC
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("E7").Style = "Normal"

# Row 8 - XRP
$ws.Range("D8").Value = "'0.621"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.94%  "
$ws.Range("E8").Style = "Normal"

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.131"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +5.69%  "
$ws.Range("E9").Style = "Normal"

# Row 10 - Cardano
$ws.Range("E10").Value = "'  -0.10%  "
$ws.Range("E10").Style = "Normal"

# Row 11 - Toncoin
$ws.Range("E11").Value = "'  -2.61%  "
$ws.Range("E11").Style = "Normal"

# Row 12 - TRON
$ws.Range("E12").Value = "'  +0.02%  "
$ws.Range("E12").Style = "Normal"

# Row 13 - Avalanche
$ws.Range("D13").Value = "'29.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.48%  "
$ws.Range("E13").Style = "Normal"

# Row 14 - ShibaInu
$ws.Range("E14").Value = "'  -1.88%  "
$ws.Range("E14").Style = "Normal"

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'3.159.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.67%  "
$ws.Range("E15").Style = "Normal"

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "'65.624.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.03%  "
$ws.Range("E16").Style = "Normal"

# Row 17 - WrappedEther
$ws.Range("D17").Value = "'2.663.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.97%  "
$ws.Range("E17").Style = "Normal"

# Row 18 - Chainlink
$ws.Range("D18").Value = "'12.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.00%  "
$ws.Range("E18").Style = "Normal"

# Row 19 - Polkadot
$ws.Range("D19").Value = "'4.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.12%  "
$ws.Range("E19").Style = "Normal"

# Row 20 - Uniswap
$ws.Range("E20").Value = "'  +1.62%  "
$ws.Range("E20").Style = "Normal"

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'352.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.92%  "
$ws.Range("E21").Style = "Normal"

# Row 22 - Dai
$ws.Range("E22").Value = "'  +0.03%  "
$ws.Range("E22").Style = "Normal"

# Row 23 - Litecoin
$ws.Range("D23").Value = "'69.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.09%  "
$ws.Range("E23").Style = "Normal"

# Row 24 - PEPE
$ws.Range("D24").Value = "'0.0000111"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +4.83%  "
$ws.Range("E24").Style = "Normal"

# Row 25 - InternetComputer(DFINITY)
$ws.Range("E25").Value = "'  -1.49%  "
$ws.Range("E25").Style = "Normal"

# Row 26 - SuiNetwork
$ws.Range("E26").Value = "'  +0.27%  "
$ws.Range("E26").Style = "Normal"

# Row 27 - Kaspa
$ws.Range("E27").Value = "'  -2.82%  "
$ws.Range("E27").Style = "Normal"

# Row 28 - Fetch.AI
$ws.Range("E28").Value = "'  -5.60%  "
$ws.Range("E28").Style = "Normal"

# Row 29 - Aptos
$ws.Range("E29").Value = "'  -3.72%  "
$ws.Range("E29").Style = "Normal"

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "'  +0.07%  "
$ws.Range("E30").Style = "Normal"

# Row 31 - Bittensor
$ws.Range("D31").Value = "'529.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.47%  "
$ws.Range("E31").Style = "Normal"

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "'  -2.38%  "
$ws.Range("E32").Style = "Normal"

# Row 33 - ImmutableX
$ws.Range("E33").Value = "'  -0.67%  "
$ws.Range("E33").Style = "Normal"

# Row 34 - RenderToken
$ws.Range("E34").Value = "'  -3.27%  "
$ws.Range("E34").Style = "Normal"

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "'  +2.39%  "
$ws.Range("E35").Style = "Normal"

# Row 36 - PolygonEcosystemToken
$ws.Range("E36").Value = "'  -1.94%  "
$ws.Range("E36").Style = "Normal"

# Row 37 - EthereumClassic
$ws.Range("D37").Value = "'20.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.32%  "
$ws.Range("E37").Style = "Normal"

# Row 38 - FirstDigitalUSD
$ws.Range("E38").Value = "'  +0.02%  "
$ws.Range("E38").Style = "Normal"

# Row 39 - Monero
$ws.Range("D39").Value = "'158.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.91%  "
$ws.Range("E39").Style = "Normal"

# Row 40 - Stacks
$ws.Range("D40").Value = "'1.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.06%  "
$ws.Range("E40").Style = "Normal"

# Row 41 - USDe (unchanged, no edits required)

# Row 42 - Aave
$ws.Range("D42").Value = "'164.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.19%  "
$ws.Range("E42").Style = "Normal"

# Row 43 - Filecoin
$ws.Range("D43").Value = "'4.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.82%  "
$ws.Range("E43").Style = "Normal"

# Row 44 - dogwifhat
$ws.Range("D44").Value = "'2.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.49%  "
$ws.Range("E44").Style = "Normal"

# Row 45 - Hedera
$ws.Range("E45").Value = "'  -0.37%  "
$ws.Range("E45").Style = "Normal"

# Row 46 - InjectiveProtocol
$ws.Range("D46").Value = "'22.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.67%  "
$ws.Range("E46").Style = "Normal"

# Row 47 - Mantle
$ws.Range("E47").Value = "'  -2.41%  "
$ws.Range("E47").Style = "Normal"

# Row 48 - BabyDogeCoin (was VeChain's slot)
$ws.Range("B48").Value = "'BabyDogeCoin"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.0₆0265"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +16.34%  "
$ws.Range("E48").Style = "Normal"

# Row 49 - VeChain (was BabyDogeCoin's slot)
$ws.Range("B49").Value = "'VeChain"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.0258"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.76%  "
$ws.Range("E49").Style = "Normal"

# Row 50 - Stellar
$ws.Range("E50").Value = "'  +1.50%  "
$ws.Range("E50").Style = "Normal"

# Row 51 - EnergySwap
$ws.Range("D51").Value = "'20.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.64%  "
$ws.Range("E51").Style = "Normal"
